# Add human-readable pipetting output
#
# Rows 16-21: sample name "CaCl2*2H2O into PHIP 1" -> "HEPES into PHIP 1",
#   PipetteConcentration (col G) updated to the new HEPES titration values.
# Row 22: becomes the "final cleaning water titration" row (was row 28).
# Row 23: becomes the "final water into water test 1" row (was row 29).
# Rows 24-29 (old duplicate HEPES + final rows) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 16-21: rename sample + update PipetteConcentration (col G) ---
$ws.Range("B16").Value = "HEPES into PHIP 1"
$ws.Range("G16").Value = 9.332155665645574

$ws.Range("B17").Value = "HEPES into PHIP 1"
$ws.Range("G17").Value = 9.832728216575232

$ws.Range("B18").Value = "HEPES into PHIP 1"
$ws.Range("G18").Value = 10.70427564894001

$ws.Range("B19").Value = "HEPES into PHIP 1"
$ws.Range("G19").Value = 9.332155665645574

$ws.Range("B20").Value = "HEPES into PHIP 1"
$ws.Range("G20").Value = 9.832728216575232

$ws.Range("B21").Value = "HEPES into PHIP 1"
$ws.Range("G21").Value = 10.70427564894001

# --- Row 22: replace with "final cleaning water titration" row content ---
$ws.Range("B22").Value = "final cleaning water titration"
$ws.Range("C22").Value = "Plates Clean.setup"
$ws.Range("D22").Value = "water5inj.inj"
$ws.Range("E22").Value = "Control"
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0

# --- Row 23: replace with "final water into water test 1" row content ---
$ws.Range("B23").Value = "final water into water test 1"
$ws.Range("D23").Value = "ChoderaWaterWater.inj"
$ws.Range("E23").Value = "Control"
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0

# --- Remove now-obsolete rows 24-29 (their content was folded into 22/23 above) ---
$ws.Range("A24:K29").EntireRow.Delete() | Out-Null
